# Automatic update of files.
# Rewrites rows 24-27 of the "Artfynd" sheet to reflect the refreshed
# export: row contents rotate (24<-27, 26<-24, 27<-26, 25 stays), Ost/Nord
# (Q/R) get rounded to whole metres, and the Starttid/Sluttid (Z/AB) time
# stamps are dropped along with a couple of stray field adjustments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 24 ----------------------------------------------------------
$ws.Range("A24").Value = 111958205
$ws.Range("B24").Value = 96348
$ws.Range("D24").Value = "VU"
$ws.Range("E24").Value = 220787
$ws.Range("F24").Value = "Knärot"
$ws.Range("G24").Value = "Goodyera repens"
$ws.Range("H24").Value = "(L.) R. Br."
# Antal (I) is stored as text in this export, not a number.
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = "1"
$ws.Range("I24").Style = "Normal"
$ws.Range("J24").Value = "plantor/tuvor"
$ws.Range("P24").Value = "Österåsen, Ång"
$ws.Range("Q24").Value = 609803
$ws.Range("R24").Value = 7011969
$ws.Range("Z24").Value = ""
$ws.Range("AB24").Value = ""
$ws.Range("AC24").Value = "½ m2"

# L24 becomes a present-but-empty cell.
$ws.Range("L24").NumberFormat = "@"
$ws.Range("L24").Value = ""
$ws.Range("L24").Style = "Normal"

# ---- Row 25 ------------------------------------------------------------
$ws.Range("Q25").Value = 609747
$ws.Range("R25").Value = 7011953
$ws.Range("Z25").Value = ""
$ws.Range("AB25").Value = ""

# ---- Row 26 --------------------------------------------------------------
$ws.Range("A26").Value = 111957843
$ws.Range("B26").Value = 89686
$ws.Range("E26").Value = 658
$ws.Range("F26").Value = "Rosenticka"
$ws.Range("G26").Value = "Rhodofomes roseus"
$ws.Range("H26").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("I26").NumberFormat = "@"
$ws.Range("I26").Value = "6"
$ws.Range("I26").Style = "Normal"
$ws.Range("J26").Value = "fruktkroppar"
$ws.Range("L26").Value = ""
$ws.Range("M26").Value = ""
$ws.Range("P26").Value = "Österås, Österås, Ång"
$ws.Range("Q26").Value = 609773
$ws.Range("R26").Value = 7011992
$ws.Range("Z26").Value = ""
$ws.Range("AB26").Value = ""
$ws.Range("AC26").Value = ""

# AF26 becomes a present-but-empty cell.
$ws.Range("AF26").NumberFormat = "@"
$ws.Range("AF26").Value = ""
$ws.Range("AF26").Style = "Normal"

# ---- Row 27 --------------------------------------------------------------
$ws.Range("A27").Value = 111958182
$ws.Range("B27").Value = 55611
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 102612
$ws.Range("F27").Value = "Järpe"
$ws.Range("G27").Value = "Tetrastes bonasia"
$ws.Range("H27").Value = "(Linnaeus, 1758)"
$ws.Range("I27").NumberFormat = "@"
$ws.Range("I27").Value = "3"
$ws.Range("I27").Style = "Normal"
$ws.Range("J27").Value = ""
$ws.Range("L27").Value = "hona"
$ws.Range("P27").Value = "Österåsen, Österås, Ång"
$ws.Range("Q27").Value = 609747
$ws.Range("R27").Value = 7011953
$ws.Range("Z27").Value = ""
$ws.Range("AB27").Value = ""
$ws.Range("AC27").Value = "1K"
$ws.Range("AF27").Value = ""

# M27 becomes a present-but-empty cell.
$ws.Range("M27").NumberFormat = "@"
$ws.Range("M27").Value = ""
$ws.Range("M27").Style = "Normal"
